$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 4 (currently 001804114 / WAGNER / 21954.32)
# for new account 005395948 / THAIS / 50000
$ws.Rows.Item(4).Insert()
$ws.Cells.Item(4, 1).NumberFormat = "@"
$ws.Cells.Item(4, 1).Value = "005395948"
$ws.Cells.Item(4, 2).Value = "THAIS"
$ws.Cells.Item(4, 3).Value = 50000

# After the insertion above, the three rows for
#   002636063 / LEDA       / 6885.31
#   003921139 / GEISA      / 6723.62
#   005231126 / WASHINGTON / 4150.15
# have shifted down to rows 8, 9, 10. Remove them and replace with a
# single new row for 005046919 / MARIANA / 2539.73.
$ws.Range("A8:C10").Delete(-4162)
$ws.Rows.Item(8).Insert()
$ws.Cells.Item(8, 1).NumberFormat = "@"
$ws.Cells.Item(8, 1).Value = "005046919"
$ws.Cells.Item(8, 2).Value = "MARIANA"
$ws.Cells.Item(8, 3).Value = 2539.73
